# Updated cryptos list (price + 1h volume/change columns) as scraped by the
# GitHub Actions job. Column D ("Price") and column E ("Volume(1h)") hold
# plain text (not numbers) in this sheet, so for price strings that look
# like a genuine number (e.g. "5.66") we prefix the literal with a leading
# apostrophe -- the standard Excel way of forcing text-entry -- so the cell
# keeps storing a string instead of being reinterpreted as a numeric value.
# Price strings that already contain multiple dots (e.g. "60.019.30") or
# other non-numeric characters don't need this, since Excel can't parse
# them as numbers anyway. The percentage strings in column E keep their
# original leading/trailing double-space padding, which likewise prevents
# Excel from auto-converting them into a numeric percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bitcoin
$ws.Range('D2').Value = '60.019.30'
$ws.Range('E2').Value = '  -0.01%  '

# Ethereum
$ws.Range('D3').Value = '2.412.82'
$ws.Range('E3').Value = '  -0.27%  '

# TetherUSD
$ws.Range('E4').Value = '  -0.01%  '

# BNB
$ws.Range('D5').Value = "'554.18"
$ws.Range('E5').Value = '  +0.42%  '

# Solana
$ws.Range('D6').Value = "'136.32"
$ws.Range('E6').Value = '  -0.69%  '

# USDC
$ws.Range('E7').Value = '  +0.03%  '

# XRP
$ws.Range('E8').Value = '  +1.09%  '

# Dogecoin
$ws.Range('E9').Value = '  -0.95%  '

# Toncoin
$ws.Range('D10').Value = "'5.66"
$ws.Range('E10').Value = '  -1.56%  '

# TRON
$ws.Range('E11').Value = '  -0.70%  '

# Cardano
$ws.Range('D12').Value = "'0.352"
$ws.Range('E12').Value = '  -1.31%  '

# Avalanche
$ws.Range('D13').Value = "'24.72"
$ws.Range('E13').Value = '  -0.40%  '

# Wrapped liquid staked Ether 2.0
$ws.Range('D14').Value = '2.845.04'
$ws.Range('E14').Value = '  -0.15%  '

# Wrapped BTC
$ws.Range('D15').Value = '59.926.65'
$ws.Range('E15').Value = '  -0.06%  '

# Shiba Inu
$ws.Range('E16').Value = '  -0.10%  '

# Wrapped Ether
$ws.Range('D17').Value = '2.402.87'
$ws.Range('E17').Value = '  -0.73%  '

# Chainlink
$ws.Range('E18').Value = '  -1.06%  '

# Bitcoin Cash
$ws.Range('D20').Value = "'326.28"
$ws.Range('E20').Value = '  -1.48%  '

# Uniswap
$ws.Range('E21').Value = '  +1.20%  '

# Dai
$ws.Range('E22').Value = '  -0.02%  '

# Litecoin
$ws.Range('D23').Value = "'64.83"
$ws.Range('E23').Value = '  -1.49%  '

# Kaspa
$ws.Range('E24').Value = '  +5.55%  '

# Internet Computer (DFINITY)
$ws.Range('D25').Value = "'8.64"
$ws.Range('E25').Value = '  -0.05%  '

# Binance-Peg BSC-USD
$ws.Range('E26').Value = '  -0.02%  '

# Fetch.AI
$ws.Range('D27').Value = "'1.40"
$ws.Range('E27').Value = '  +3.85%  '

# PancakeSwap
$ws.Range('E28').Value = '  -0.04%  '

# PEPE
$ws.Range('D29').Value = '0.0₃0771'
$ws.Range('E29').Value = '  -1.24%  '

# Monero
$ws.Range('D30').Value = "'170.85"
$ws.Range('E30').Value = '  +0.22%  '

# Aptos
$ws.Range('E31').Value = '  -1.66%  '

# Sui Network
$ws.Range('D32').Value = "'1.09"
$ws.Range('E32').Value = '  +6.34%  '

# Polygon Ecosystem Token
$ws.Range('D33').Value = "'0.400"
$ws.Range('E33').Value = '  -3.26%  '

# Ethereum Classic
$ws.Range('D34').Value = "'18.38"
$ws.Range('E34').Value = '  -1.26%  '

# USDe
$ws.Range('E35').Value = '  +0.02%  '

# ImmutableX
$ws.Range('E36').Value = '  +2.12%  '

# FirstDigitalUSD
$ws.Range('E37').Value = '  +0.03%  '

# NEAR Protocol
$ws.Range('D38').Value = "'4.21"
$ws.Range('E38').Value = '  +0.89%  '

# Bittensor
$ws.Range('D39').Value = "'325.04"
$ws.Range('E39').Value = '  +3.37%  '

# Stacks
$ws.Range('D40').Value = "'1.59"
$ws.Range('E40').Value = '  -1.35%  '

# Aave
$ws.Range('D41').Value = "'145.97"
$ws.Range('E41').Value = '  +5.10%  '

# Filecoin
$ws.Range('D42').Value = "'3.60"
$ws.Range('E42').Value = '  -1.99%  '

# Stellar
$ws.Range('E43').Value = '  +0.13%  '

# Injective Protocol
$ws.Range('D44').Value = "'19.79"
$ws.Range('E44').Value = '  +2.29%  '

# Hedera
$ws.Range('E45').Value = '  -0.88%  '

# Mantle
$ws.Range('D46').Value = "'0.577"
$ws.Range('E46').Value = '  -0.25%  '

# VeChain
$ws.Range('E47').Value = '  -1.70%  '

# WhiteBIT Coin
$ws.Range('D48').Value = "'11.04"
$ws.Range('E48').Value = '  -0.01%  '

# dogwifhat
$ws.Range('D49').Value = "'1.57"
$ws.Range('E49').Value = '  -1.55%  '

# ZEEBU
$ws.Range('E50').Value = '  -0.62%  '

# BitgetToken
$ws.Range('D51').Value = "'0.938"
$ws.Range('E51').Value = '  -1.80%  '
